$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly report rows above the current row 155, pushing the
# existing rows 155:160 down to 157:162 (values/styles travel with them).
$ws.Rows("155:156").Insert()

# Row 155 - new weekly entry (Primera, $/caja 18 kilos)
$ws.Range("A155").Value = 9
$ws.Range("B155").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C155").Value = "Metropolitana"
$ws.Range("D155").Value = 45075
$ws.Range("E155").Value = 13
$ws.Range("F155").Value = 100114002
$ws.Range("G155").Value = "Camote"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 790
$ws.Range("K155").Value = 18000
$ws.Range("L155").Value = 19000
$ws.Range("M155").Value = 18494
$ws.Range("N155").Value = "$/caja 18 kilos"
$ws.Range("O155").Value = "Perú"
$ws.Range("P155").Value = 1027
$ws.Range("Q155").Value = 18
$ws.Range("R155").Value = "Hortaliza"

# Row 156 - new weekly entry (Primera, $/malla 18 kilos)
$ws.Range("A156").Value = 9
$ws.Range("B156").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C156").Value = "Metropolitana"
$ws.Range("D156").Value = 45075
$ws.Range("E156").Value = 13
$ws.Range("F156").Value = 100114002
$ws.Range("G156").Value = "Camote"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 880
$ws.Range("K156").Value = 15000
$ws.Range("L156").Value = 16000
$ws.Range("M156").Value = 15500
$ws.Range("N156").Value = "$/malla 18 kilos"
$ws.Range("O156").Value = "Perú"
$ws.Range("P156").Value = 861
$ws.Range("Q156").Value = 18
$ws.Range("R156").Value = "Hortaliza"
